# feat: implement funtionality to download excel
#
# Rebuilds the "Worksheet" sheet with the full set of client/service
# records that are exported when the report is downloaded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1, 1).Value = "Nombre"
$ws.Cells.Item(1, 2).Value = "Teléfono"
$ws.Cells.Item(1, 3).Value = "Servicio"

# Data rows: Nombre, Teléfono, Servicio
$rows = @(
    @("fiorella", "bustamante", "cámaras"),
    @(3434, 344343343, "internet"),
    @(3434, 34434334, "internet"),
    @(3434, 3, "internet"),
    @("Esteban Salas Sulca", 927859435, "cámaras"),
    @("Esteban Salas Sulca", 927859435, "cámaras"),
    @("Esteban Salas Sulca", 927859435, "internet"),
    @("Esteban Salas Sulca", 927859435, "cable"),
    @("Esteban Salas Sulca", 927859435, "cable"),
    @("Esteban Salas Sulca", 927859435, "cable"),
    @("Esteban Salas Sulca", 927859435, "internet"),
    @("ge", 927859435, "internet_cable"),
    @("Esteban Salas Sulca", 927859435, "internet"),
    @("Esteban Salas Sulca", 927859435, "internet"),
    @("Esteban Salas Sulca", 927859435, "internet_cable"),
    @("Esteban Salas Sulca", 927859435, "internet")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
